$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Column B (the topic "Guide" slug column) shifts up by one slot
#    because the "ci" entry was dropped from the underlying topic
#    list. Rewrite B4:B59 (the block affected by the shift); B2/B3
#    keep their existing text but are re-written too for safety.
# ------------------------------------------------------------------
$slugs = @(
  "accessibility",
  "browser-webpage-performance",
  "code-editors",
  "command-line-interface",
  "css-fonts-and-icons",
  "css-fundamentals",
  "css-in-js",
  "css-layout",
  "css-media-queries",
  "css-practices",
  "css-tooling",
  "css-transitions-transforms-&-animations",
  "css-ui-toolkits",
  "data-apis",
  "desktop-apps",
  "dom-bom-cssom",
  "front-end-development-tech-overview",
  "getting-a-front-end-developer-job",
  "git",
  "html",
  "html-email",
  "http",
  "index.html",
  "jam-stack",
  "javascript-fundamentals",
  "javascript-modules-scripts",
  "javascript-regular-expressions",
  "js-animation",
  "js-app-frameworks",
  "js-async",
  "js-bundlers",
  "js-compile-to-js",
  "js-cs",
  "js-fp",
  "js-oop",
  "js-performance",
  "js-practices-tools",
  "js-runtime-ast",
  "js-state-management",
  "js-ui-components-and-widgets",
  "js-web-api",
  "mobile-apps",
  "node-npm",
  "npm-yarn-scripts",
  "pwa",
  "rwd",
  "seo",
  "template",
  "testing",
  "the-front-end-developer-profession",
  "the-web-platform",
  "the-www-basic-mechanics",
  "ui-design-patterns",
  "web-apps",
  "web-images",
  "web-security",
  "webdev-tools",
  "x-compile-to-js"
)

for ($i = 0; $i -lt $slugs.Length; $i++) {
    $ws.Cells.Item(2 + $i, 2).Value = $slugs[$i]
}

# ------------------------------------------------------------------
# 2. Column C ("Topic" / generated-title source) gets real values for
#    the first 8 data rows.
# ------------------------------------------------------------------
$ws.Range("C2").Value = "Learn Accessibility"
$ws.Range("C3").Value = "Learn Browser, Website, and Web App Performance"
$ws.Range("C4").Value = "Learn Code Editors"
$ws.Range("C5").Value = "Learn the Command Line Interface"
$ws.Range("C6").Value = "Learn CSS Fonts & Icons"
$ws.Range("C7").Value = "Learn CSS Fundamentals"
$ws.Range("C8").Value = "Learn CSS in JavaScript"
$ws.Range("C9").Value = "Learn CSS Layout"

# ------------------------------------------------------------------
# 3. Column D ("Description") gets a test value in row 2.
# ------------------------------------------------------------------
$ws.Range("D2").Value = "Test description"

# ------------------------------------------------------------------
# 4. Header row: new column F header + re-affirm existing headers.
# ------------------------------------------------------------------
$ws.Range("B1").Value = "Guide"
$ws.Range("C1").Value = "Topic"
$ws.Range("D1").Value = "Description"
$ws.Range("E1").Value = "Generated Title Tag"
$ws.Range("F1").Value = "Generated Meta Tag"
$ws.Range("F1").Style = $ws.Range("E1").Style

# ------------------------------------------------------------------
# 5. Column F: new "Generated Meta Tag" formula column, mirroring the
#    existing E-column shared-formula pattern.
# ------------------------------------------------------------------
$ws.Range("F2:F60").Formula = "=""<meta name='description' content=''""&D2&""'>"""

# ------------------------------------------------------------------
# 6. Remove row 60 (its topic "x-compile-to-js" now lives at row 59
#    after the column-B shift, so the former last row is dropped).
# ------------------------------------------------------------------
$ws.Rows.Item(60).Delete()

# ------------------------------------------------------------------
# 7. Re-affirm the E-column shared title formula across the new
#    (shorter) range so it stays a single shared formula block.
# ------------------------------------------------------------------
$ws.Range("E3:E59").Formula = "=""<title>Learning ""&C3&"" Resources - Front-End Developer Learning Roadmap</title>"""

# ------------------------------------------------------------------
# 8. Column widths.
# ------------------------------------------------------------------
$ws.Columns.Item(3).ColumnWidth = 55.1640625
$ws.Columns.Item(5).ColumnWidth = 104.83203125
$ws.Columns.Item(6).ColumnWidth = 85.5

# ------------------------------------------------------------------
# 9. Selection.
# ------------------------------------------------------------------
$ws.Range("C9").Select()
